# [Fonds de solidarite] Add 2020-08-11 data
# Updates nombre_aides (col C) and montant_total (col D) for the rows whose
# figures were refreshed with the 2020-08-11 data pull. Values are written
# with a leading apostrophe so Excel stores them as text (matching the
# workbook's existing inline-string / text-typed columns) and the style is
# then reset to "Normal" so no stray number-format / quote-prefix styling
# is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Count = "422";  Amount = "981828.79" },
    @{ Row = 4;  Count = "856";  Amount = "2655616.22" },
    @{ Row = 6;  Count = "542";  Amount = "1566513.59" },
    @{ Row = 28; Count = "253";  Amount = "646542.64" },
    @{ Row = 30; Count = "512";  Amount = "2030800.70" },
    @{ Row = 32; Count = "360";  Amount = "1175478.17" },
    @{ Row = 33; Count = "11";   Amount = "34500.00" },
    @{ Row = 45; Count = "322";  Amount = "885357.74" },
    @{ Row = 47; Count = "565";  Amount = "2151830.99" },
    @{ Row = 48; Count = "379";  Amount = "1262567.16" },
    @{ Row = 51; Count = "3335"; Amount = "7600616.99" },
    @{ Row = 52; Count = "23";   Amount = "142000.00" },
    @{ Row = 53; Count = "3804"; Amount = "12889689.79" },
    @{ Row = 55; Count = "3887"; Amount = "11901642.08" },
    @{ Row = 61; Count = "130";  Amount = "287500.00" },
    @{ Row = 62; Count = "222";  Amount = "603400.00" },
    @{ Row = 65; Count = "19";   Amount = "74110.00" },
    @{ Row = 66; Count = "18";   Amount = "68496.05" },
    @{ Row = 67; Count = "5";    Amount = "27700.00" },
    @{ Row = 68; Count = "227";  Amount = "597331.00" },
    @{ Row = 69; Count = "367";  Amount = "1233844.26" },
    @{ Row = 70; Count = "207";  Amount = "649916.51" },
    @{ Row = 72; Count = "17";   Amount = "56000.00" },
    @{ Row = 73; Count = "366";  Amount = "899135.70" },
    @{ Row = 75; Count = "881";  Amount = "2942919.89" },
    @{ Row = 76; Count = "498";  Amount = "1606208.87" },
    @{ Row = 77; Count = "35";   Amount = "94000.00" },
    @{ Row = 78; Count = "33";   Amount = "133736.09" }
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)   # column C: nombre_aides
    $cCell.Value = "'" + $u.Count
    $cCell.Style = "Normal"

    $dCell = $ws.Cells.Item($u.Row, 4)   # column D: montant_total
    $dCell.Value = "'" + $u.Amount
    $dCell.Style = "Normal"
}
